# Fruta / hortaliza, semanal
#
# Inserts one new weekly price record for "Feria Lagunitas de Puerto Montt"
# (Mandarina) ahead of the existing row 193, pushing the historical rows
# 193-243 down to 194-244 (dimension grows from A1:T243 to A1:T244).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 193 - everything currently at/after row 193
# (including its formatting) shifts down by one row.
$ws.Rows(193).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(193, 1).Value  = 4
$ws.Cells.Item(193, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(193, 3).Value  = "Los Lagos"
$ws.Cells.Item(193, 4).Value  = 44754
$ws.Cells.Item(193, 5).Value  = 10
$ws.Cells.Item(193, 6).Value  = "Fruta"
$ws.Cells.Item(193, 7).Value  = 100102
$ws.Cells.Item(193, 8).Value  = "Cítricos"
$ws.Cells.Item(193, 9).Value  = 100102004
$ws.Cells.Item(193, 10).Value = "Mandarina"
$ws.Cells.Item(193, 11).Value = "Clemenuless"
$ws.Cells.Item(193, 12).Value = "Primera"
$ws.Cells.Item(193, 13).Value = 800
$ws.Cells.Item(193, 14).Value = 9000
$ws.Cells.Item(193, 15).Value = 10000
$ws.Cells.Item(193, 16).Value = 9500
$ws.Cells.Item(193, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(193, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(193, 19).Value = 950
$ws.Cells.Item(193, 20).Value = 10
